$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a pure number and would otherwise be
# auto-converted to a numeric cell by Excel; pre-format these as Text so the
# literal string (e.g. leading/trailing zeros) is preserved exactly.
$textCells = @(
    "D5",
    "D6",
    "D8",
    "D11",
    "D16",
    "D20",
    "D21",
    "D24",
    "D27",
    "D28",
    "D31",
    "D33",
    "D34",
    "D37",
    "D39",
    "D41",
    "D43",
    "D45",
    "D47",
    "D50",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply every cell-value change from the diff.
$ws.Range("D2").Value = "60.794.55"
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").Value = "2.907.14"
$ws.Range("E3").Value = "  -2.73%  "
$ws.Range("D5").Value = "526.17"
$ws.Range("E5").Value = "  -2.87%  "
$ws.Range("D6").Value = "144.30"
$ws.Range("E6").Value = "  -4.71%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.546"
$ws.Range("D9").Value = "2.914.70"
$ws.Range("E9").Value = "  -2.89%  "
$ws.Range("E10").Value = "  -5.06%  "
$ws.Range("D11").Value = "6.15"
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("E12").Value = "  -2.90%  "
$ws.Range("D13").Value = "3.415.08"
$ws.Range("E13").Value = "  -2.73%  "
$ws.Range("E14").Value = "  +2.55%  "
$ws.Range("D15").Value = "60.820.04"
$ws.Range("E15").Value = "  -1.59%  "
$ws.Range("D16").Value = "22.54"
$ws.Range("E16").Value = "  -5.97%  "
$ws.Range("D17").Value = "2.909.32"
$ws.Range("E17").Value = "  -2.94%  "
$ws.Range("E18").Value = "  -3.81%  "
$ws.Range("E19").Value = "  -5.11%  "
$ws.Range("D20").Value = "11.55"
$ws.Range("E20").Value = "  -4.20%  "
$ws.Range("D21").Value = "353.33"
$ws.Range("E21").Value = "  -6.41%  "
$ws.Range("E22").Value = "  -2.81%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "5.66"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  -1.94%  "
$ws.Range("E26").Value = "  -3.94%  "
$ws.Range("D27").Value = "0.178"
$ws.Range("E27").Value = "  -5.00%  "
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("E29").Value = "  -6.95%  "
$ws.Range("E30").Value = "  -5.13%  "
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("E32").Value = "  -2.42%  "
$ws.Range("D33").Value = "19.64"
$ws.Range("E33").Value = "  -3.85%  "
$ws.Range("D34").Value = "153.27"
$ws.Range("E34").Value = "  -4.77%  "
$ws.Range("E35").Value = "  -4.17%  "
$ws.Range("E36").Value = "  -6.43%  "
$ws.Range("D37").Value = "0.996"
$ws.Range("E37").Value = "  -6.66%  "
$ws.Range("E38").Value = "  -6.16%  "
$ws.Range("D39").Value = "37.56"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("E40").Value = "  -5.43%  "
$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").Value = "0.653"
$ws.Range("E41").Value = "  -2.83%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.291.19"
$ws.Range("E42").Value = "  -5.13%  "
$ws.Range("D43").Value = "3.70"
$ws.Range("E43").Value = "  -5.18%  "
$ws.Range("E44").Value = "  -1.21%  "
$ws.Range("D45").Value = "20.37"
$ws.Range("E45").Value = "  -7.63%  "
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").Value = "4.93"
$ws.Range("E47").Value = "  -4.69%  "
$ws.Range("E48").Value = "  -3.15%  "
$ws.Range("E49").Value = "  -0.64%  "
$ws.Range("D50").Value = "0.0913"
$ws.Range("D51").Value = "18.54"
$ws.Range("E51").Value = "  -5.66%  "
